$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_G"

$ws.Range("A2").Value = 52.5460930640913
$ws.Range("A3").Value = 52.5460930640913
$ws.Range("A4").Value = 52.589991220368745
$ws.Range("A5").Value = 49.60491659350307
$ws.Range("A6").Value = 50.043898156277436
$ws.Range("A7").Value = 51.00965759438104
$ws.Range("A8").Value = 54.302019315188765
$ws.Range("A9").Value = 55.04828797190518
$ws.Range("A10").Value = 54.38981562774363
$ws.Range("A11").Value = 54.43371378402107
$ws.Range("A12").Value = 59.35030728709394
$ws.Range("A13").Value = 59.1747146619842
$ws.Range("A14").Value = 53.819139596136964
$ws.Range("A15").Value = 54.69710272168569
$ws.Range("A16").Value = 54.08252853380158
$ws.Range("A17").Value = 53.46795434591747
$ws.Range("A18").Value = 54.87269534679543
$ws.Range("A19").Value = 54.74100087796313
$ws.Range("A20").Value = 53.424056189640034
$ws.Range("A21").Value = 53.204565408252854
$ws.Range("A22").Value = 53.55575065847235
$ws.Range("A23").Value = 56.277436347673394
$ws.Range("A24").Value = 61.15013169446883
$ws.Range("A25").Value = 60.7550482879719
$ws.Range("A26").Value = 53.950834064969264
$ws.Range("A27").Value = 53.950834064969264
$ws.Range("A28").Value = 54.25812115891132
$ws.Range("A29").Value = 61.062335381913954
$ws.Range("A30").Value = 60.447761194029844
$ws.Range("A31").Value = 54.03863037752414
$ws.Range("A32").Value = 56.36523266022827
$ws.Range("A33").Value = 52.01931518876207
$ws.Range("A34").Value = 52.23880597014925
$ws.Range("A35").Value = 55.355575065847226
$ws.Range("A36").Value = 55.70676031606673
$ws.Range("A37").Value = 58.691834942932395
$ws.Range("A38").Value = 56.18964003511853
$ws.Range("A39").Value = 55.79455662862159
$ws.Range("A40").Value = 56.8920105355575
$ws.Range("A41").Value = 55.00438981562774
$ws.Range("A42").Value = 54.78489903424057
$ws.Range("A43").Value = 55.17998244073748
$ws.Range("A44").Value = 53.8630377524144
$ws.Range("A45").Value = 54.03863037752414
$ws.Range("A46").Value = 53.204565408252854
$ws.Range("A47").Value = 53.73134328358209
$ws.Range("A48").Value = 52.282704126426694
$ws.Range("A49").Value = 53.77524143985952
$ws.Range("A50").Value = 50.39508340649693
